# Weekly fruit/vegetable price update: a new "Apio" (celery) price record was
# reported for the week, so insert a new row at position 100 (pushing the
# existing rows 100-110 down to 101-111) and populate it with the new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("100:100").Insert()

$ws.Range("A100").Value = 7
$ws.Range("B100").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C100").Value = "Ñuble"
$ws.Range("D100").Value = 44449
$ws.Range("E100").Value = 16
$ws.Range("F100").Value = 100112017
$ws.Range("G100").Value = "Apio"
$ws.Range("H100").Value = "Americana (o)"
$ws.Range("I100").Value = "Primera"
$ws.Range("J100").Value = 160
$ws.Range("K100").Value = 8000
$ws.Range("L100").Value = 8500
$ws.Range("M100").Value = 8250
$ws.Range("N100").Value = "$/docena de matas"
$ws.Range("O100").Value = "Provincia del Elquí"
$ws.Range("P100").Value = 1375
$ws.Range("Q100").Value = 6
$ws.Range("R100").Value = "Hortaliza"
